$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "C1IDc8d501"
$ws.Range("B1").Value = "1:1, 2:2, 3:2, 3:3, 4:3, 4:4, 3:4, 2:4, 1:4, 1:3"
$ws.Range("C1").Value = "Shining, Bles brothers, Bles brothers, 12 hommes en coleres, 12 hommes en coleres, Hellraiser, Hellraiser, Hellraiser, Hellraiser, 12 hommes en coleres"
$ws.Range("D1").Value = 'Affiche, R\xe9alisater, Ann\xe9e, Ann\xe9e, Genre, Genre, Ann\xe9e, R\xe9alisater, Affiche, Affiche'
$ws.Range("E1").Value = "12 hommes en coleres"
$ws.Range("F1").Value = "12 hommes en coleres"
